$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (id, repositoryUrl, repositoryName, repositoryAuthor, startingDate, OSE, BCE, PDE, SV, OS, SD, RS, TFS, UI, TC)
$rows = @(
    @(172, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0"),
    @(173, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0"),
    @(174, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0"),
    @(175, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0"),
    @(176, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0"),
    @(177, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "1", "1", "1", "1", "0", "0", "1", "1"),
    @(178, "https://github.com/tensorflow/ranking", "ranking", "tensorflow", "12/03/2018", "0", "0", "0", "1", "1", "1", "0", "0", "1", "0"),
    @(179, "https://github.com/rescrv/HyperDex", "HyperDex", "rescrv", "04/13/2011", "1", "0", "1", "1", "0", "0", "0", "0", "0", "1"),
    @(180, "https://github.com/pocl/pocl", "pocl", "pocl", "02/08/2011", "0", "0", "1", "1", "1", "0", "0", "0", "1", "1"),
    @(181, "https://github.com/raghakot/keras-vis", "keras-vis", "raghakot", "11/11/2016", "0", "1", "1", "1", "1", "1", "0", "0", "0", "1")
)

$startRow = 173
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A: numeric id, styled like the rest of the id column.
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item(172, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    # Columns B-O: stored as literal text (urls, dates-as-text, "0"/"1" flags),
    # matching the sheet's existing inline-string convention. The leading
    # apostrophe forces text entry (avoids Excel's date/number auto-detect);
    # re-applying the base cell style clears the quote-prefix formatting
    # that the apostrophe entry leaves behind.
    for ($c = 2; $c -le $data.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = "'" + $data[$c - 1]
        $cell.Style = "Normale"
    }
}
